$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 6 with new values (previously blank A6/B6, C6/D6 didn't exist)
$ws.Range("A6").Value = "GFG/CN"
$ws.Range("B6").Value = "GFG/CN"
$ws.Range("C6").Value = "0 - 1 Knapsack Problem"
$ws.Range("D6").Value = "Java"

# Make sure C6 and D6 pick up the same formatting as the rest of row 6 (A6/B6):
# left-aligned, top-aligned, wrapped text (matches existing style of A6/B6)
$ws.Range("C6:D6").HorizontalAlignment = -4131 # xlLeft
$ws.Range("C6:D6").VerticalAlignment = -4160   # xlTop
$ws.Range("C6:D6").WrapText = $true

# Update the active selection to D6 (was D10)
$ws.Range("D6").Select()
